$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = -0.0  # A2
$ws.Cells.Item(2, 2).Value = -0.0736143064681548  # B2
$ws.Cells.Item(2, 3).Value = -0.0  # C2
$ws.Cells.Item(2, 4).Value = 0.2023283625086515  # D2
$ws.Cells.Item(2, 5).Value = 0.005361901149070607  # E2
$ws.Cells.Item(2, 7).Value = 0.0  # G2
$ws.Cells.Item(2, 9).Value = -0.0  # I2
$ws.Cells.Item(2, 10).Value = -0.0  # J2
$ws.Cells.Item(2, 11).Value = 0.005910374655943606  # K2
$ws.Cells.Item(2, 12).Value = -0.0  # L2
$ws.Cells.Item(2, 13).Value = 0.2017405004068997  # M2
$ws.Cells.Item(2, 14).Value = -0.003175673222564392  # N2
$ws.Cells.Item(2, 18).Value = -0.0  # R2
$ws.Cells.Item(2, 19).Value = 0.0  # S2
$ws.Cells.Item(2, 20).Value = -0.0836551046379089  # T2
$ws.Cells.Item(2, 22).Value = 0.01416007321150566  # V2
$ws.Cells.Item(2, 23).Value = -0.03298341659304817  # W2
$ws.Cells.Item(2, 25).Value = -0.0  # Y2
$ws.Cells.Item(2, 26).Value = -0.0  # Z2
$ws.Cells.Item(2, 28).Value = 0.0  # AB2
$ws.Cells.Item(2, 29).Value = -0.05416417955287071  # AC2
$ws.Cells.Item(2, 30).Value = 0.0  # AD2
$ws.Cells.Item(2, 31).Value = -0.01788162495550331  # AE2
$ws.Cells.Item(2, 32).Value = 0.0002069792777307436  # AF2
$ws.Cells.Item(2, 33).Value = -0.0  # AG2
$ws.Cells.Item(2, 34).Value = -0.0  # AH2
$ws.Cells.Item(2, 35).Value = -0.0  # AI2
$ws.Cells.Item(2, 36).Value = 0.0  # AJ2
$ws.Cells.Item(2, 37).Value = -0.0  # AK2
$ws.Cells.Item(2, 38).Value = -0.03184932082569965  # AL2
$ws.Cells.Item(2, 39).Value = 0.0  # AM2
$ws.Cells.Item(2, 40).Value = 0.02827770634814052  # AN2
$ws.Cells.Item(2, 41).Value = 0.06944358562979185  # AO2
$ws.Cells.Item(2, 43).Value = 0.0  # AQ2
$ws.Cells.Item(2, 44).Value = -0.0  # AR2
$ws.Cells.Item(2, 46).Value = 0.0  # AT2
$ws.Cells.Item(2, 47).Value = -0.1497027310705481  # AU2
$ws.Cells.Item(2, 49).Value = 0.07080831603100772  # AW2
$ws.Cells.Item(2, 50).Value = -0.001817919973327277  # AX2
$ws.Cells.Item(2, 51).Value = -0.0  # AY2
$ws.Cells.Item(2, 55).Value = -0.0  # BC2
$ws.Cells.Item(2, 56).Value = -0.01418805710578807  # BD2
$ws.Cells.Item(2, 58).Value = 0.08649163433815991  # BF2
$ws.Cells.Item(2, 59).Value = 0.03283123518905573  # BG2
$ws.Cells.Item(2, 62).Value = -0.0  # BJ2
$ws.Cells.Item(2, 64).Value = 0.0  # BL2
$ws.Cells.Item(2, 65).Value = 0.03190481457958391  # BM2
$ws.Cells.Item(2, 67).Value = -0.04257541708426302  # BO2
$ws.Cells.Item(2, 68).Value = -0.08905310676590357  # BP2
$ws.Cells.Item(2, 73).Value = 0.0  # BU2
$ws.Cells.Item(2, 74).Value = -0.04640710802875297  # BV2
$ws.Cells.Item(2, 75).Value = 0.0  # BW2
$ws.Cells.Item(2, 76).Value = 0.01094127294829059  # BX2
$ws.Cells.Item(2, 77).Value = -0.02004983166574863  # BY2
$ws.Cells.Item(2, 78).Value = -0.0  # BZ2
$ws.Cells.Item(2, 80).Value = 0.0  # CB2
$ws.Cells.Item(2, 82).Value = -0.0  # CD2
$ws.Cells.Item(2, 83).Value = 0.03274941356648393  # CE2
$ws.Cells.Item(2, 85).Value = -0.03099117283538349  # CG2
$ws.Cells.Item(2, 86).Value = 0.01593166186343111  # CH2
$ws.Cells.Item(2, 88).Value = -0.0  # CJ2
$ws.Cells.Item(2, 91).Value = -0.0  # CM2
$ws.Cells.Item(2, 92).Value = -0.01061381960660221  # CN2
$ws.Cells.Item(2, 94).Value = 0.02133540246658532  # CP2
$ws.Cells.Item(2, 95).Value = 0.03716094318380431  # CQ2
$ws.Cells.Item(2, 98).Value = 0.0  # CT2
$ws.Cells.Item(2, 99).Value = -0.0  # CU2
$ws.Cells.Item(2, 100).Value = -0.0  # CV2
$ws.Cells.Item(2, 101).Value = 0.04597314766486385  # CW2
$ws.Cells.Item(2, 103).Value = -0.03342648399499332  # CY2
$ws.Cells.Item(2, 104).Value = 0.01017704690408558  # CZ2
$ws.Cells.Item(2, 108).Value = -0.0  # DD2
$ws.Cells.Item(2, 109).Value = -0.0  # DE2
$ws.Cells.Item(2, 110).Value = 0.02874569132567836  # DF2
$ws.Cells.Item(2, 112).Value = 0.02900526664094873  # DH2
$ws.Cells.Item(2, 113).Value = 0.03380050877759293  # DI2
$ws.Cells.Item(2, 114).Value = 0.0  # DJ2
$ws.Cells.Item(2, 115).Value = -0.0  # DK2
$ws.Cells.Item(2, 116).Value = -0.0  # DL2
$ws.Cells.Item(2, 118).Value = 0.0  # DN2
$ws.Cells.Item(2, 119).Value = -0.01950247745448723  # DO2
$ws.Cells.Item(2, 120).Value = -0.0  # DP2
$ws.Cells.Item(2, 121).Value = 0.03637034262361485  # DQ2
$ws.Cells.Item(2, 122).Value = -0.01945341551444906  # DR2
$ws.Cells.Item(2, 123).Value = -0.0  # DS2
$ws.Cells.Item(2, 127).Value = 0.0  # DW2
$ws.Cells.Item(2, 128).Value = -0.05668836815106189  # DX2
$ws.Cells.Item(2, 129).Value = -0.0  # DY2
$ws.Cells.Item(2, 130).Value = -0.008248668484950638  # DZ2
$ws.Cells.Item(2, 131).Value = -0.02434073422596091  # EA2
$ws.Cells.Item(2, 132).Value = 0.0  # EB2
$ws.Cells.Item(2, 136).Value = -0.0  # EF2
$ws.Cells.Item(2, 137).Value = 0.04085971144248264  # EG2
$ws.Cells.Item(2, 139).Value = 0.06686095049629477  # EI2
$ws.Cells.Item(2, 140).Value = -0.02407360759003618  # EJ2
$ws.Cells.Item(2, 145).Value = 0.0  # EO2
$ws.Cells.Item(2, 146).Value = 0.04575396385905522  # EP2
$ws.Cells.Item(2, 147).Value = 0.0  # EQ2
$ws.Cells.Item(2, 148).Value = -0.0349068518066118  # ER2
$ws.Cells.Item(2, 149).Value = 0.03654119765287879  # ES2
$ws.Cells.Item(2, 150).Value = 0.0  # ET2
$ws.Cells.Item(2, 151).Value = -0.0  # EU2
$ws.Cells.Item(2, 152).Value = 0.0  # EV2
$ws.Cells.Item(2, 154).Value = 0.0  # EX2
$ws.Cells.Item(2, 155).Value = 0.04340351386436194  # EY2
$ws.Cells.Item(2, 157).Value = -0.02669773779825179  # FA2
$ws.Cells.Item(2, 158).Value = 0.01710822948871973  # FB2
$ws.Cells.Item(2, 160).Value = -0.0  # FD2
$ws.Cells.Item(2, 163).Value = -0.0  # FG2
$ws.Cells.Item(2, 164).Value = 0.001691558334483329  # FH2
$ws.Cells.Item(2, 165).Value = 0.0  # FI2
$ws.Cells.Item(2, 166).Value = -0.006288705109680439  # FJ2
$ws.Cells.Item(2, 167).Value = -0.008476023910300627  # FK2
$ws.Cells.Item(2, 168).Value = -0.0  # FL2
$ws.Cells.Item(2, 170).Value = -0.0  # FN2
$ws.Cells.Item(2, 172).Value = -0.0  # FP2
$ws.Cells.Item(2, 173).Value = -0.01406165917757605  # FQ2
$ws.Cells.Item(2, 174).Value = -0.0  # FR2
$ws.Cells.Item(2, 175).Value = -0.0184191401413391  # FS2
$ws.Cells.Item(2, 176).Value = 0.005010644615690384  # FT2
$ws.Cells.Item(2, 178).Value = -0.0  # FV2
$ws.Cells.Item(2, 179).Value = -0.0  # FW2
$ws.Cells.Item(2, 181).Value = 0.0  # FY2
$ws.Cells.Item(2, 182).Value = -0.03040407898268115  # FZ2
$ws.Cells.Item(2, 184).Value = 0.03399166782388464  # GB2
$ws.Cells.Item(2, 186).Value = 0.0  # GD2
$ws.Cells.Item(2, 187).Value = -0.0  # GE2

Write-Host "Updated row 2 values"